$d = $word.ActiveDocument

# The document has two distinct logo pictures, each placed twice (once per
# header/footer variant - "primary" and "first page"):
#   - a BTec logo picture living in both headers, currently named "image2.jpg",
#     being renamed to "image1.jpg"
#   - a Pearson logo picture living in both footers, currently named "image1.png",
#     being renamed to "image2.png"
#
# Renaming is done through InlineShape.Name, reached via the Selection object
# (Selection.InlineShapes) rather than Range.InlineShapes directly - for the
# footer pictures specifically, assigning .Name on the Range-addressed
# InlineShape silently fails to stick, while going through an explicit
# Select() + Selection.InlineShapes(1).Name = ... reliably commits the change.

$section = $d.Sections(1)

# --- Headers: BTec logo, image2.jpg -> image1.jpg ---
for ($i = 1; $i -le 2; $i++) {
    $hdr = $section.Headers($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes(1)
        [void]$shp.Select()
        $word.Selection.InlineShapes(1).Name = "image1.jpg"
        Write-Output ("Header " + $i + ": renamed logo to image1.jpg")
    }
}

# --- Footers: Pearson logo, image1.png -> image2.png ---
for ($i = 1; $i -le 2; $i++) {
    $ftr = $section.Footers($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes(1)
        [void]$shp.Select()
        $word.Selection.InlineShapes(1).Name = "image2.png"
        Write-Output ("Footer " + $i + ": renamed logo to image2.png")
    }
}
